$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AZ1").Value = 0.9282672939276414
$ws.Range("BC1").Value = 0.89313939542188803
$ws.Range("BP1").Value = 0.94121117777457197
$ws.Range("A2").Value = 0.99800389874655604
$ws.Range("D2").Value = 0.89747708800335468
$ws.Range("A3").Value = 0.91956940906341222
$ws.Range("D3").Value = 0.85558679932231008
$ws.Range("F4").Value = 0.99914076759258619
$ws.Range("AI4").Value = 0.58230124119830595
$ws.Range("C5").Value = 0.70090550600508728
$ws.Range("G5").Value = 0.89346422619533672
$ws.Range("F7").Value = 0.78215171857753041
$ws.Range("H7").Value = 0.98915500464602002
$ws.Range("V7").Value = 0.99014877737842566
$ws.Range("AQ7").Value = 0.79698906394168301
$ws.Range("BL7").Value = 0.62425620828102635
$ws.Range("I8").Value = 0.57582147103865178
$ws.Range("H10").Value = 0.97442870805484194
$ws.Range("I10").Value = 0.81370635863887131
$ws.Range("BI10").Value = 0.67230629378271312
$ws.Range("F11").Value = 0.59135786781686517
$ws.Range("I11").Value = 0.72826355286717925
$ws.Range("M12").Value = 0.93946845063648299
$ws.Range("N12").Value = 0.98768764147763077
$ws.Range("S12").Value = 0.94696709222551245
$ws.Range("N13").Value = 0.97278893960390089
$ws.Range("O13").Value = 0.5891077961554233
$ws.Range("AH13").Value = 0.80572600767189728
$ws.Range("AF14").Value = 0.98096953672310272
$ws.Range("BN14").Value = 0.81550976530504493
$ws.Range("N15").Value = 0.92903291863993187
$ws.Range("P15").Value = 0.92043240866822296
$ws.Range("AW15").Value = 0.91139951938420571
$ws.Range("Q16").Value = 0.78327874060529323
$ws.Range("AM16").Value = 0.84053789440823001
$ws.Range("O17").Value = 0.82867412729060907
$ws.Range("BP17").Value = 0.83265491746081421
$ws.Range("Q18").Value = 0.58744350487491248
$ws.Range("S18").Value = 0.86159561540376584
$ws.Range("BN18").Value = 0.99329832181440836
$ws.Range("U19").Value = 0.88718304305517748
$ws.Range("T21").Value = 0.98247597010592225
$ws.Range("T22").Value = 0.87034863654915007
$ws.Range("U22").Value = 0.65430581391238751
$ws.Range("W22").Value = 0.84307581280897748
$ws.Range("BG22").Value = 0.84607543474482427
$ws.Range("U23").Value = 0.72366701116183196
$ws.Range("AF23").Value = 0.60317504778448594
$ws.Range("W25").Value = 0.65106751142542085
$ws.Range("AA25").Value = 0.87622574375404461
$ws.Range("X26").Value = 0.83547931422793864
$ws.Range("Y26").Value = 0.82602411252759311
$ws.Range("AB26").Value = 0.88783556149878085
$ws.Range("Z27").Value = 0.79643681909913044
$ws.Range("AC27").Value = 0.76977543738494325
$ws.Range("AW27").Value = 0.85803169814841551
$ws.Range("AC28").Value = 0.96868095636096907
$ws.Range("AO28").Value = 0.93126653365118595
$ws.Range("AC30").Value = 0.71895950118733443
$ws.Range("AE30").Value = 0.87630119236813964
$ws.Range("AF30").Value = 0.91942690341689082
$ws.Range("S31").Value = 0.79026508947704288
$ws.Range("AC31").Value = 0.96284703719937992
$ws.Range("B32").Value = 0.73869846457744526
$ws.Range("AG32").Value = 0.64492639066729396
$ws.Range("AR32").Value = 0.73266437508713711
$ws.Range("AE33").Value = 0.88854474572399489
$ws.Range("AZ33").Value = 0.78698135928509583
$ws.Range("BC34").Value = 0.91624845697740875
$ws.Range("AG35").Value = 0.93969857976078974
$ws.Range("AU35").Value = 0.68177284170478025
$ws.Range("AA36").Value = 0.71242730293447054
$ws.Range("AH36").Value = 0.7854452579590151
$ws.Range("K37").Value = 0.87541473987192786
$ws.Range("AI37").Value = 0.94779997319838682
$ws.Range("AM37").Value = 0.7938940162496384
$ws.Range("J38").Value = 0.93763259821355138
$ws.Range("AJ38").Value = 0.91019130103063461
$ws.Range("AK38").Value = 0.98484071412053331
$ws.Range("AM38").Value = 0.69432201520450798
$ws.Range("AN38").Value = 0.91966266733704805
$ws.Range("AN39").Value = 0.82657182229935011
$ws.Range("AO40").Value = 0.92888262525691179
$ws.Range("AQ41").Value = 0.95929510127961448
$ws.Range("AN42").Value = 0.9288715518210422
$ws.Range("AO42").Value = 0.92762746238120997
$ws.Range("AQ42").Value = 0.83303737433683045
$ws.Range("AP44").Value = 0.79387891212352302
$ws.Range("AT44").Value = 0.63380223421014581
$ws.Range("AQ45").Value = 0.63017718182791682
$ws.Range("AT45").Value = 0.79177284795520064
$ws.Range("BB45").Value = 0.93359786164109138
$ws.Range("BO45").Value = 0.74401422389922489
$ws.Range("AS47").Value = 0.97086342085046184
$ws.Range("AT47").Value = 0.92412859733278918
$ws.Range("AW47").Value = 0.7822818612671032
$ws.Range("AT48").Value = 0.68482623487574812
$ws.Range("AU48").Value = 0.86358813657389533
$ws.Range("AV49").Value = 0.70298966372946092
$ws.Range("AX49").Value = 0.98933409445660181
$ws.Range("AF50").Value = 0.95451366387932279
$ws.Range("AU50").Value = 0.65798118262280503
$ws.Range("E51").Value = 0.90517501014504387
$ws.Range("AZ51").Value = 0.95724999405877575
$ws.Range("BA51").Value = 0.92047835839830494
$ws.Range("BB52").Value = 0.95808087247843687
$ws.Range("BB53").Value = 0.74911093706547338
$ws.Range("AV54").Value = 0.87151011558997626
$ws.Range("AX54").Value = 0.85410221171185607
$ws.Range("BA55").Value = 0.98793559689048671
$ws.Range("BB55").Value = 0.82852914186642712
$ws.Range("BE55").Value = 0.92999950034804568
$ws.Range("BF56").Value = 0.92970672461887571
$ws.Range("K57").Value = 0.62391189494184185
$ws.Range("X57").Value = 0.74469370597889528
$ws.Range("BD57").Value = 0.73348319270579276
$ws.Range("BE58").Value = 0.92584770316919851
$ws.Range("BG58").Value = 0.97239419852927
$ws.Range("BH58").Value = 0.89788049810998327
$ws.Range("AW59").Value = 0.85010431634660422
$ws.Range("BI59").Value = 0.98573330283580107
$ws.Range("E60").Value = 0.55867823165764396
$ws.Range("BB60").Value = 0.96197854012218453
$ws.Range("BJ61").Value = 0.84484338302922635
$ws.Range("BK61").Value = 0.98351971278973571
$ws.Range("BL62").Value = 0.95908537270905792
$ws.Range("BJ63").Value = 0.95992265239625885
$ws.Range("BL63").Value = 0.9683051922340149
$ws.Range("BM63").Value = 0.70224094171370499
$ws.Range("AM64").Value = 0.82510945154178428
$ws.Range("BN65").Value = 0.92243021816937776
$ws.Range("BH66").Value = 0.96883800190332359
$ws.Range("BL66").Value = 0.85328204507452421
$ws.Range("BM67").Value = 0.76140502669641807
$ws.Range("BP67").Value = 0.99078382575789536
$ws.Range("B68").Value = 0.63781044040291657
